# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Ragnarok_Profits workbook tabs (one per crafting profession).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1377.9714
$ws.Range("I15").Value = 1377.9714
$ws.Range("K15").Value = 4133.914199999999
$ws.Range("M15").Value = -3964.914199999999
$ws.Range("H40").Value = 250001950
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 250001950
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 250001950
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -250002300
$ws.Range("H107").Value = 837.0769
$ws.Range("J107").Value = 1455.125
$ws.Range("L107").Value = 1455.125
$ws.Range("N107").Value = -5295.125
$ws.Range("H113").Value = 8619.25
$ws.Range("I113").Value = 8619.25
$ws.Range("K113").Value = 8619.25
$ws.Range("M113").Value = -5365.25
$ws.Range("H137").Value = 2109.4375
$ws.Range("I137").Value = 2087.9285
$ws.Range("K137").Value = 6263.7855
$ws.Range("M137").Value = -3713.7855
$ws.Range("H138").Value = 6323.72
$ws.Range("I138").Value = 2278.4666
$ws.Range("J138").Value = 12391.6
$ws.Range("K138").Value = 6835.399800000001
$ws.Range("L138").Value = 37174.8
$ws.Range("M138").Value = -1695.399800000001
$ws.Range("N138").Value = -47454.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3121
$ws.Range("I2").Value = 3083.0715
$ws.Range("K2").Value = 3083.0715
$ws.Range("M2").Value = -2970.0715
$ws.Range("H32").Value = 3127.806
$ws.Range("I32").Value = 2278.3447
$ws.Range("J32").Value = 8602.111000000001
$ws.Range("K32").Value = 2278.3447
$ws.Range("L32").Value = 8602.111000000001
$ws.Range("M32").Value = -1991.3447
$ws.Range("N32").Value = -9176.111000000001
$ws.Range("H61").Value = 17694018
$ws.Range("I61").Value = 21001824
$ws.Range("K61").Value = 21001824
$ws.Range("M61").Value = -21001612
$ws.Range("H110").Value = 2947.625
$ws.Range("I110").Value = 896.25
$ws.Range("J110").Value = 4999
$ws.Range("K110").Value = 896.25
$ws.Range("L110").Value = 4999
$ws.Range("M110").Value = 1148.75
$ws.Range("N110").Value = -9089
$ws.Range("H116").Value = 3121
$ws.Range("I116").Value = 3083.0715
$ws.Range("K116").Value = 3083.0715
$ws.Range("M116").Value = -789.0715
$ws.Range("H132").Value = 2277589
$ws.Range("I132").Value = 4392.75
$ws.Range("K132").Value = 13178.25
$ws.Range("M132").Value = -10648.25
$ws.Range("H136").Value = 17694018
$ws.Range("I136").Value = 21001824
$ws.Range("K136").Value = 63005472
$ws.Range("M136").Value = -63002922

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3121
$ws.Range("I3").Value = 3083.0715
$ws.Range("K3").Value = 3083.0715
$ws.Range("M3").Value = -2969.0715
$ws.Range("H134").Value = 6669538.5
$ws.Range("I134").Value = 2923.9167
$ws.Range("J134").Value = 33335996
$ws.Range("K134").Value = 8771.750100000001
$ws.Range("L134").Value = 100007988
$ws.Range("M134").Value = -6236.750100000001
$ws.Range("N134").Value = -100013058

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5897146.5
$ws.Range("I16").Value = 6265530.5
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 6265530.5
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -6265243.5
$ws.Range("N16").Value = -3574
$ws.Range("H31").Value = 43482150
$ws.Range("J31").Value = 3670.3635
$ws.Range("L31").Value = 3670.3635
$ws.Range("N31").Value = -4260.363499999999
$ws.Range("H34").Value = 43482150
$ws.Range("J34").Value = 3670.3635
$ws.Range("L34").Value = 3670.3635
$ws.Range("N34").Value = -4074.3635
$ws.Range("H58").Value = 2783.6
$ws.Range("J58").Value = 3014
$ws.Range("L58").Value = 3014
$ws.Range("N58").Value = -3420
$ws.Range("H105").Value = 1961
$ws.Range("I105").Value = 1402.8125
$ws.Range("K105").Value = 1402.8125
$ws.Range("M105").Value = 344.1875
$ws.Range("H113").Value = 5897146.5
$ws.Range("I113").Value = 6265530.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 6265530.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -6263360.5
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 3994.9
$ws.Range("I122").Value = 3334.8333
$ws.Range("K122").Value = 10004.4999
$ws.Range("M122").Value = -7554.499899999999
$ws.Range("H132").Value = 2825.423
$ws.Range("I132").Value = 2764.8948
$ws.Range("J132").Value = 2989.7144
$ws.Range("K132").Value = 8294.6844
$ws.Range("L132").Value = 8969.143199999999
$ws.Range("M132").Value = -5764.6844
$ws.Range("N132").Value = -14029.1432
$ws.Range("H136").Value = 2783.6
$ws.Range("J136").Value = 3014
$ws.Range("L136").Value = 9042
$ws.Range("N136").Value = -14142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J80").Value = 21635
$ws.Range("L80").Value = 64905
$ws.Range("N80").Value = -66777
$ws.Range("J83").Value = 21635
$ws.Range("L83").Value = 194715
$ws.Range("N83").Value = -204075
$ws.Range("H107").Value = 6501161
$ws.Range("I107").Value = 357.66666
$ws.Range("J107").Value = 8274107
$ws.Range("K107").Value = 1072.99998
$ws.Range("L107").Value = 24822321
$ws.Range("M107").Value = 847.0000199999999
$ws.Range("N107").Value = -24826161
$ws.Range("H113").Value = 913.6667
$ws.Range("I113").Value = 733.25
$ws.Range("K113").Value = 2199.75
$ws.Range("M113").Value = -29.75
$ws.Range("H134").Value = 22245.521
$ws.Range("I134").Value = 2045.125
$ws.Range("J134").Value = 33019.066
$ws.Range("K134").Value = 6135.375
$ws.Range("L134").Value = 99057.198
$ws.Range("M134").Value = -1065.375
$ws.Range("N134").Value = -109197.198

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 903.53845
$ws.Range("I97").Value = 1045.2
$ws.Range("K97").Value = 1045.2
$ws.Range("M97").Value = -549.2
$ws.Range("H102").Value = 3180.1052
$ws.Range("I102").Value = 3235.5334
$ws.Range("K102").Value = 3235.5334
$ws.Range("M102").Value = -1613.5334
$ws.Range("H113").Value = 1325728.2
$ws.Range("I113").Value = 3199.1428
$ws.Range("J113").Value = 2648257.2
$ws.Range("K113").Value = 3199.1428
$ws.Range("L113").Value = 2648257.2
$ws.Range("M113").Value = -1029.1428
$ws.Range("N113").Value = -2652597.2
$ws.Range("H132").Value = 7697174
$ws.Range("I132").Value = 4595.778
$ws.Range("K132").Value = 13787.334
$ws.Range("M132").Value = -11257.334
$ws.Range("H134").Value = 89569.25
$ws.Range("J134").Value = 89569.25
$ws.Range("L134").Value = 268707.75
$ws.Range("N134").Value = -273777.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2184.8
$ws.Range("J46").Value = 2641.6667
$ws.Range("L46").Value = 2641.6667
$ws.Range("N46").Value = -3017.6667
$ws.Range("H55").Value = 1065.0358
$ws.Range("I55").Value = 527.13336
$ws.Range("J55").Value = 1685.6923
$ws.Range("K55").Value = 527.13336
$ws.Range("L55").Value = 1685.6923
$ws.Range("M55").Value = -354.13336
$ws.Range("N55").Value = -2031.6923
$ws.Range("H122").Value = 3596.8604
$ws.Range("I122").Value = 3296.8157
$ws.Range("K122").Value = 9890.447100000001
$ws.Range("M122").Value = -7440.447100000001
$ws.Range("H132").Value = 5197.154
$ws.Range("I132").Value = 3313.2
$ws.Range("J132").Value = 6374.625
$ws.Range("K132").Value = 9939.599999999999
$ws.Range("L132").Value = 19123.875
$ws.Range("M132").Value = -7409.599999999999
$ws.Range("N132").Value = -24183.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 671.7273
$ws.Range("I113").Value = 568.0833
$ws.Range("K113").Value = 1704.2499
$ws.Range("M113").Value = 465.7501
$ws.Range("H132").Value = 263162.06
$ws.Range("I132").Value = 6770.485
$ws.Range("J132").Value = 1673315.6
$ws.Range("K132").Value = 20311.455
$ws.Range("L132").Value = 5019946.800000001
$ws.Range("M132").Value = -17781.455
$ws.Range("N132").Value = -5025006.800000001
